$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Force the Price column to be stored as text so values like "1.00" or "240.92"
# are not silently converted to numbers by Excel, matching the source data which
# is plain text throughout column D.
$ws.Range("D2:D51").NumberFormat = "@"

$ws.Range("D2").Value = "96.548.24"
$ws.Range("E2").Value = "  -0.69%  "
$ws.Range("D3").Value = "3.675.63"
$ws.Range("E3").Value = "  +1.96%  "
$ws.Range("E4").Value = "  +0.10%  "
$ws.Range("D5").Value = "240.92"
$ws.Range("E5").Value = "  -1.09%  "
$ws.Range("D6").Value = "1.84"
$ws.Range("E6").Value = "  +8.99%  "
$ws.Range("D7").Value = "663.82"
$ws.Range("E7").Value = "  +0.65%  "
$ws.Range("D8").Value = "0.422"
$ws.Range("E8").Value = "  +0.57%  "
$ws.Range("E9").Value = "  +1.43%  "
$ws.Range("E10").Value = "  +0.07%  "
$ws.Range("D11").Value = "3.671.63"
$ws.Range("E11").Value = "  +1.84%  "
$ws.Range("D12").Value = "45.30"
$ws.Range("E12").Value = "  +3.04%  "
$ws.Range("E13").Value = "  +0.32%  "
$ws.Range("D14").Value = "6.85"
$ws.Range("E14").Value = "  +5.77%  "
$ws.Range("D15").Value = "4.360.66"
$ws.Range("E15").Value = "  +2.02%  "
$ws.Range("D16").Value = "0.0000268"
$ws.Range("E16").Value = "  +2.70%  "
$ws.Range("D17").Value = "96.270.84"
$ws.Range("E17").Value = "  -0.65%  "
$ws.Range("D18").Value = "8.85"
$ws.Range("E18").Value = "  +14.03%  "
$ws.Range("D19").Value = "3.665.85"
$ws.Range("E19").Value = "  +1.70%  "
$ws.Range("D20").Value = "12.84"
$ws.Range("E20").Value = "  +0.91%  "
$ws.Range("D21").Value = "18.46"
$ws.Range("E21").Value = "  +2.09%  "
$ws.Range("D22").Value = "0.523"
$ws.Range("E22").Value = "  -2.18%  "
$ws.Range("D23").Value = "526.14"
$ws.Range("E23").Value = "  +2.57%  "
$ws.Range("D24").Value = "3.43"
$ws.Range("E24").Value = "  -0.06%  "
$ws.Range("E25").Value = "  -1.30%  "
$ws.Range("D26").Value = "6.98"
$ws.Range("E26").Value = "  +1.32%  "
$ws.Range("D27").Value = "101.84"
$ws.Range("E27").Value = "  +3.46%  "
$ws.Range("D28").Value = "13.04"
$ws.Range("E28").Value = "  -0.48%  "
$ws.Range("B29").Value = "WrappedeETH"
$ws.Range("C29").Value = "https://coinranking.com/coin/dFlrSolOX+wrappedeeth-weeth"
$ws.Range("D29").Value = "3.875.68"
$ws.Range("E29").Value = "  +2.13%  "
$ws.Range("B30").Value = "Hedera"
$ws.Range("C30").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D30").Value = "0.169"
$ws.Range("E30").Value = "  +12.03%  "
$ws.Range("B31").Value = "InternetComputer(DFINITY)"
$ws.Range("C31").Value = "https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp"
$ws.Range("D31").Value = "12.48"
$ws.Range("E31").Value = "  +6.29%  "
$ws.Range("B32").Value = "PancakeSwap"
$ws.Range("C32").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D32").Value = "3.06"
$ws.Range("E32").Value = "  +0.15%  "
$ws.Range("B33").Value = "Dai"
$ws.Range("C33").Value = "https://coinranking.com/coin/MoTuySvg7+dai-dai"
$ws.Range("D33").Value = "0.999"
$ws.Range("E33").Value = "  -0.05%  "
$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.91"
$ws.Range("E34").Value = "  +16.88%  "
$ws.Range("B35").Value = "Cronos"
$ws.Range("C35").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D35").Value = "0.186"
$ws.Range("E35").Value = "  -0.51%  "
$ws.Range("B36").Value = "EthereumClassic"
$ws.Range("C36").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D36").Value = "32.75"
$ws.Range("E36").Value = "  +3.04%  "
$ws.Range("B37").Value = "Binance-PegBSC-USD"
$ws.Range("C37").Value = "https://coinranking.com/coin/i5jggxiwp+binance-pegbsc-usd-bsc-usd"
$ws.Range("D37").Value = "1.00"
$ws.Range("E37").Value = "  -0.45%  "
$ws.Range("B38").Value = "Bittensor"
$ws.Range("C38").Value = "https://coinranking.com/coin/pgv7xSFi6+bittensor-tao"
$ws.Range("D38").Value = "644.03"
$ws.Range("E38").Value = "  +3.67%  "
$ws.Range("B39").Value = "PolygonEcosystemToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol"
$ws.Range("D39").Value = "0.593"
$ws.Range("E39").Value = "  +3.46%  "
$ws.Range("B40").Value = "RenderToken"
$ws.Range("C40").Value = "https://coinranking.com/coin/vfo5XYwcV+rendertoken-render"
$ws.Range("D40").Value = "8.76"
$ws.Range("E40").Value = "  -0.02%  "
$ws.Range("B41").Value = "EnergySwap"
$ws.Range("C41").Value = "https://coinranking.com/coin/SbWqqTui-+energyswap-ens"
$ws.Range("D41").Value = "44.46"
$ws.Range("E41").Value = "  +32.87%  "
$ws.Range("B42").Value = "Kaspa"
$ws.Range("C42").Value = "https://coinranking.com/coin/V8GxkwWow+kaspa-kas"
$ws.Range("D42").Value = "0.161"
$ws.Range("E42").Value = "  +4.74%  "
$ws.Range("B43").Value = "ARBITRUM"
$ws.Range("C43").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D43").Value = "0.966"
$ws.Range("E43").Value = "  +4.12%  "
$ws.Range("B44").Value = "ImmutableX"
$ws.Range("C44").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D44").Value = "1.97"
$ws.Range("E44").Value = "  +2.15%  "
$ws.Range("B45").Value = "Filecoin"
$ws.Range("C45").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D45").Value = "6.40"
$ws.Range("E45").Value = "  +7.98%  "
$ws.Range("B46").Value = "USDe"
$ws.Range("C46").Value = "https://coinranking.com/coin/exbfr2U-0+usde-usde"
$ws.Range("D46").Value = "1.00"
$ws.Range("E46").Value = "  -0.01%  "
$ws.Range("B47").Value = "Algorand"
$ws.Range("C47").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D47").Value = "0.455"
$ws.Range("E47").Value = "  +22.75%  "
$ws.Range("B48").Value = "VeChain"
$ws.Range("C48").Value = "https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet"
$ws.Range("D48").Value = "0.0463"
$ws.Range("E48").Value = "  +6.58%  "
$ws.Range("B49").Value = "Stacks"
$ws.Range("C49").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D49").Value = "2.30"
$ws.Range("E49").Value = "  -1.63%  "
$ws.Range("B50").Value = "WhiteBITCoin"
$ws.Range("C50").Value = "https://coinranking.com/coin/GE4c3_TbB+whitebitcoin-wbt"
$ws.Range("D50").Value = "23.65"
$ws.Range("E50").Value = "  -0.16%  "
$ws.Range("B51").Value = "MantraDAO"
$ws.Range("C51").Value = "https://coinranking.com/coin/cTdD8lD-6+mantradao-om"
$ws.Range("D51").Value = "3.66"
$ws.Range("E51").Value = "  +3.73%  "
